# Insert a new data row at row 64 (pushing existing rows 64-108 down to 65-109)
# and populate it with the new "Granada" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 64, shifting rows 64:108 down to 65:109
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new record's data
$ws.Range("A64").Value = 10
$ws.Range("B64").Value = "Vega Modelo de Temuco"
$ws.Range("C64").Value = "La Araucanía"
$ws.Range("D64").Value = 44658
$ws.Range("E64").Value = 9
$ws.Range("F64").Value = "Fruta"
$ws.Range("G64").Value = 100104
$ws.Range("H64").Value = "Frutos de pepita"
$ws.Range("I64").Value = 100104001
$ws.Range("J64").Value = "Granada"
$ws.Range("K64").Value = "Wonderfull"
$ws.Range("L64").Value = "Primera"
$ws.Range("M64").Value = 200
$ws.Range("N64").Value = 11000
$ws.Range("O64").Value = 12000
$ws.Range("P64").Value = 11500
$ws.Range("Q64").Value = "$/bandeja 10 kilos granel"
$ws.Range("R64").Value = "Provincia de Limarí"
$ws.Range("S64").Value = 1150
$ws.Range("T64").Value = 10
